# Refactored and QoL changes
# Turn the filled-in "Harsh Morayya" sample row into a blank header/label
# template row (Roll no. / Name / Batch / Specialization / email / Semester /
# phone no. / Program code) plus a placeholder second row telling the user
# where to start entering data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header labels -----------------------------------------------
# (E1 is handled further down, after the hyperlink is rebuilt, so the
#  "email" text isn't clobbered by Hyperlinks.Add's TextToDisplay arg.)
$ws.Range("A1").Value = "Roll no."
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Batch"
$ws.Range("D1").Value = "Specialization"

# --- E1: rebuild the mailto hyperlink with an explicit display text ----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E1"), "mailto:hars@sas.com", "", "", "hars@sas.com")
$ws.Range("E1").Style = "Hyperlink"
$ws.Range("E1").Value = "email"

# --- Remaining row 1 header labels (F1/G1 were numbers before) ---------
$ws.Range("F1").Value = "Semester"
$ws.Range("G1").Value = "phone no."
$ws.Range("H1").Value = "Program code"

# --- Row 2: placeholder prompting the user to fill in real data --------
$ws.Range("A2").Value = "<Start inserting data in this row>"

# --- Selection moves to A2, ready for data entry ------------------------
[void]$ws.Range("A2").Select()
